$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F holds "想去人数" (want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 2088
$ws1.Range("F6").Value = 638
$ws1.Range("F8").Value = 2076
$ws1.Range("F9").Value = 10722
$ws1.Range("F11").Value = 158
$ws1.Range("F15").Value = 7576
$ws1.Range("F17").Value = 722
$ws1.Range("F18").Value = 269
$ws1.Range("F20").Value = 3342

# Sheet "全部类型" (All Types) - same events, column F counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 2088
$ws4.Range("F6").Value = 638
$ws4.Range("F9").Value = 2076
$ws4.Range("F12").Value = 10722
$ws4.Range("F14").Value = 158
$ws4.Range("F18").Value = 7576
$ws4.Range("F20").Value = 722
$ws4.Range("F21").Value = 269
$ws4.Range("F23").Value = 3342
